$wb = $excel.ActiveWorkbook

$oldStr = "January 30 2026 16.19.47 EST"
$newStr = "February 02 2026 12.49.33 EST"

# "About" sheet: update version text in A2 and the citation text in A6
$wsAbout = $wb.Worksheets.Item("About")

$a2 = $wsAbout.Cells.Item(2, 1).Value2
$wsAbout.Cells.Item(2, 1).Value2 = $a2.Replace($oldStr, $newStr)

$a6 = $wsAbout.Cells.Item(6, 1).Value2
$wsAbout.Cells.Item(6, 1).Value2 = $a6.Replace($oldStr, $newStr)

# "Boundaries and methane sources" sheet: update build_version column (S) for rows 2-7
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

foreach ($r in 2..7) {
    $cell = $wsData.Cells.Item($r, 19)  # Column S = 19
    $val = $cell.Value2
    $cell.Value2 = $val.Replace($oldStr, $newStr)
}
